$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Color Data set delete" / winner-declaration counter -------------------
# The num1/num2/num3 counter block (rows 3-14) used to count UP from 2 to 13
# (A4 = A3+1, filled down). Flip it so the block counts DOWN from 13 to 2
# instead (A4 = A3-1, filled down) - this is the "delete" direction used by
# the winner-declaration function.

# New seed values for row 3.
$ws.Range("A3").Value = 13
$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 13

# Row 4: change the operator from +1 to -1.
$ws.Range("A4").Formula = "=A3-1"
$ws.Range("B4").Formula = "=B3-1"
$ws.Range("C4").Formula = "=C3-1"

# Rows 5-14: same "-1 from the cell above" pattern, written explicitly per
# row so relative references advance one row at a time (equivalent to
# dragging the fill handle from row 4 down to row 14).
for ($r = 5; $r -le 14; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 1).Formula = "=A$prev-1"
    $ws.Cells.Item($r, 2).Formula = "=B$prev-1"
    $ws.Cells.Item($r, 3).Formula = "=C$prev-1"
}

# Further down, the num3 (column C) block that starts its own countdown at
# C18 now starts one lower (13 -> 12); the C19:C26 formulas recompute off of
# that automatically.
$ws.Range("C18").Value = 12

# --- View state ---------------------------------------------------------
# Zoom back in to 98% (was 52%), drop the scrolled-down top-left cell, and
# move the active selection up to F7 (was L451).
$excel.ActiveWindow.Zoom = 98
$ws.Range("F7").Select()
